# "New Excel sheets provided by Tom (and minor bug fixes)"
#
# Bug fixes: two stray/incorrect Infrastructure-ID values ("I-8" in E2 and
# "I-1" in E9) are cleared out on the Infrastructure Connections column.
# E2 becomes a blank placeholder (single space, matching the sheet's usual
# "blank" convention), and E9 is cleared to a true empty cell. Once those
# shared strings are no longer referenced anywhere, Excel's save drops the
# now-unused "I-8"/"I-1" entries from the shared-strings table, shifting
# every later string's index down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = " "
$ws.Range("E9").Value = ""

# Leave the cursor where the author left it when they saved.
$null = $ws.Range("D10").Select()
